$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded. It belongs chronologically
# right after the header block of this product's entries, so insert a
# fresh row at row 15 (pushing the existing rows 15-24 down to 16-25,
# matching the diff's <dimension ref="A1:T25"/>).
$ws.Rows("15:15").Insert()

# The new row mirrors the row that is now directly below it (the old
# row 15, now row 16) for every column except the date - same market,
# region, product taxonomy, volume, prices, unit and origin.
$ws.Range("A16:T16").Copy()
$ws.Range("A15").PasteSpecial()
$excel.CutCopyMode = $false

# Only the date (column D) differs for the new observation.
$ws.Range("D15").Value = 44629
